# Update outdated URL on About tab
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Replace the old EPA URL text with the new one and attach a hyperlink to it.
$ws.Hyperlinks.Add(
    $ws.Range("B6"),
    "https://www.epa.gov/environmental-economics/mortality-risk-valuation",
    "whatvalue",
    "",
    "https://www.epa.gov/environmental-economics/mortality-risk-valuation - whatvalue"
)
$ws.Range("B6").Value = "https://www.epa.gov/environmental-economics/mortality-risk-valuation#whatvalue"

# Adding the hyperlink re-applies the Hyperlink cell style via a fresh style
# record; restore the original "Hyperlink" style assignment on B6.
$ws.Range("B6").Style = "Hyperlink"
